$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing 20 data rows (A2:C21) down by one row, to A3:C22,
# to make room for the newly inserted row at row 2.
$existing = $ws.Range("A2:C21").Value2
$ws.Range("A3:C22").Value2 = $existing

# Write the newly inserted row's values at row 2.
$ws.Range("A2").Value2 = -0.772541880607605
$ws.Range("B2").Value2 = -0.2669219076633453
$ws.Range("C2").Value2 = 0.974086582660675

# Append the 9 new rows of data at the bottom (rows 23-31).
$newRows = @(
    @(-0.4122457504272466, -0.09671294689178367, 0.6507992744445799),
    @(-0.01871716976165649, 0.1319747567176817, -0.1633049249649066),
    @(-0.1193938255310074, -0.05726575851440492, -0.4602591991424554),
    @(-0.2904316186904892, -0.1059370636939995, -0.1864967942237846),
    @(0.0775488615036011, -0.0401190519332886, -0.02606511116027833),
    @(-0.02678942680358923, -0.3534234166145316, -0.111013770103454),
    @(-0.1420207023620607, -0.1315011978149412, 0.04810285568237234),
    @(-0.1982678174972537, -0.0857929587364194, -0.1000801920890816),
    @(0.08944976329803658, 0.1226030588150037, -0.6240378618240389)
)

$row = 23
foreach ($values in $newRows) {
    $ws.Cells.Item($row, 1).Value2 = $values[0]
    $ws.Cells.Item($row, 2).Value2 = $values[1]
    $ws.Cells.Item($row, 3).Value2 = $values[2]
    $row++
}
